$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Property1"
$ws2 = $wb.Worksheets.Item(2)   # "Record_Cooldown"

# --- Fix formatting on row 6 (B6 / G6 previously had mismatched styles that
#     forced numeric values to be stored as text). Copy the plain numeric
#     style used by the rest of the row (C6) onto B6 and G6 before writing
#     the new values so they stay numeric. ---
$ws1.Range("C6").Copy()
$ws1.Range("B6").PasteSpecial(-4122)   # xlPasteFormats
$ws1.Range("G6").PasteSpecial(-4122)   # xlPasteFormats

# --- Cell value updates ---
$ws1.Range("G3").Value = 0

$ws1.Range("B6").Value = 1
$ws1.Range("C6").Value = 1
$ws1.Range("D6").Value = 1
$ws1.Range("E6").Value = 1
$ws1.Range("F6").Value = 1
$ws1.Range("G6").Value = 1

# --- Extend the TRUE/FALSE list validation to also cover B6:E6 ---
$ws1.Range("B6:E6").Validation.Add(3, 1, 1, "TRUE,FALSE")

# --- Selection / active sheet changes: Property1 becomes the active tab with
#     G3 selected (moving away from Record_Cooldown which loses tabSelected) ---
$ws1.Activate()
$ws1.Range("G3").Select()
